$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 into the new I1:J1 header cells, then set labels
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I0 and IF values for data rows 2..77
$iVals = @(9,9,9,9,9,10,9,8,9,9,9,7,10,8,9,9,9,9,9,9,9,9,8,8,8,9,10,8,9,8,7,8,8,8,8,9,9,9,9,9,8,8,9,8,9,8,8,8,8,9,6,9,9,8,7,9,8,9,9,9,9,9,9,10,9,8,9,9,9,9,6,5,8,7,5,5)
$jVals = @(9,9,9,9,9,10,9,9,9,9,9,7,10,8,9,9,9,9,9,9,9,9,9,8,8,9,11,8,9,8,7,8,8,9,9,9,9,9,10,9,9,8,9,8,9,9,8,8,8,9,6,9,9,8,7,9,8,9,9,9,9,9,9,11,9,8,9,9,9,9,6,5,8,7,5,5)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
